$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.779.67', '  +0.23%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.599.91', '  +0.04%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '''1.01', '  +0.14%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '''211.57', '  +0.06%  '),
    @(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '''0.512', '  -0.06%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '''1.00', '  +0.17%  '),
    @(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '''0.0620', '  +0.20%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '''0.247', '  -0.14%  '),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '''19.72', '  +1.13%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '''0.0848', '  +0.93%  '),
    @(12, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.822.35', '  -0.10%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.590.35', '  +2.27%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '''4.07', '  +0.69%  '),
    @(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '''0.526', '  +0.52%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '''65.17', '  -0.25%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '26.705.26', '  +0.06%  '),
    @(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₃0743', '  -2.09%  '),
    @(19, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '''209.85', '  +0.05%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '''1.01', '  +0.23%  '),
    @(21, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '''7.14', '  -0.18%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '''4.30', '  +0.54%  '),
    @(23, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '''2.25', '  -2.05%  '),
    @(24, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '''9.03', '  +1.03%  '),
    @(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '''144.16', '  +0.72%  '),
    @(26, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '''1.01', '  +0.03%  '),
    @(27, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '''7.14', '  +0.06%  '),
    @(28, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '''0.114', '  -0.46%  '),
    @(29, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '''15.37', '  +0.34%  '),
    @(30, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '''0.0509', '  -1.93%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '''1.15', '  -0.48%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '''3.27', '  +0.74%  '),
    @(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '''2.99', '  +1.05%  '),
    @(34, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '''1.27', '  +17.63%  '),
    @(35, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.278.35', '  -0.82%  '),
    @(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '''2.49', '  +1.05%  '),
    @(37, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '''1.49', '  -0.23%  '),
    @(38, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '''0.597', '  -3.41%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '''0.0169', '  -1.40%  '),
    @(40, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '''0.827', '  +0.02%  '),
    @(41, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '''2.24', '  +2.38%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '''5.46', '  +0.26%  '),
    @(43, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '''0.778', '  -0.82%  '),
    @(44, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '''62.77', '  -0.48%  '),
    @(45, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.734.12', '  -0.17%  '),
    @(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '''90.48', '  -0.94%  '),
    @(47, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '''1.57', '  -0.21%  '),
    @(48, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '''0.103', '  +2.44%  '),
    @(49, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '''0.0512', '  +0.89%  '),
    @(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '''7.55', '  +2.44%  '),
    @(51, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.0₇0965', '  -7.51%  '),
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}
